$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Purchase 22-23") ---
# Remove the two leading invoice rows of the "Collective Trade Links Pvt Ltd" group
# (I-C-1-23-445419 / I-C-1-23-445420) and the next two rows
# (I-C-1-23-450387 / I-C-1-23-450796), leaving only the last two invoices of that
# group (I-C-1-23-450927 / I-C-1-23-451040).
# Also remove the whole "Digiserve" group and the whole "V M Traders" group.
# Delete bottom-up so row numbers used in each step stay valid.
$ws1.Range("A30:A31").EntireRow.Delete()   # V M Traders group (entirely removed)
$ws1.Range("A25:A26").EntireRow.Delete()   # Digiserve row + following blank spacer row
$ws1.Range("A18:A21").EntireRow.Delete()   # first four rows of the Collective Trade Links group

# Fix up the sequence numbers / totals that Excel cannot infer automatically.
# The Collective Trade Links group now starts at row 18, keep its group number "2".
$ws1.Range("A18").Value = 2
# Rebuild the now-orphaned subtotal formula for that group (only 2 data rows remain).
$ws1.Range("F19").Formula = "=E18+E19"
# The old "Pilz India Pvt Ltd" group (previously group 4) is now group 3.
$ws1.Range("A21").Value = 3

$ws1.Range("A24:XFD26").Select()

# --- Sheet2 ("Sale 22-23") ---
$ws2.Range("A26").Select()

$ws1.Activate()
